$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AW2").Value = 170.939699
$ws.Range("AW3").Value = 14.055694
$ws.Range("AW4").Value = 1.813738
$ws.Range("AW5").Value = 125.885359
$ws.Range("AW6").Value = 128.67191
$ws.Range("AK7").Value = 28.777025
$ws.Range("AW8").Value = 163.70941
$ws.Range("AW9").Value = 112.750046
$ws.Range("AW10").Value = 96.833333
$ws.Range("AW11").Value = 12.729977
$ws.Range("AW12").Value = 97.735706
$ws.Range("AW13").Value = 175.793403
$ws.Range("AW14").Value = 161.608553
$ws.Range("AW15").Value = 130.737083
$ws.Range("AW16").Value = 133.034387
$ws.Range("AW17").Value = 12.684201
$ws.Range("AQ18").Value = 18.640891
$ws.Range("AW19").Value = 100.621215
$ws.Range("AW20").Value = 92.709977
$ws.Range("AW21").Value = 14.056123
$ws.Range("AK22").Value = 28.777037
$ws.Range("AK23").Value = 28.777083
$ws.Range("AW24").Value = 125.884826
$ws.Range("AW25").Value = 12.829039
$ws.Range("AK26").Value = 28.777188
$ws.Range("AW27").Value = 161.614063
$ws.Range("AW28").Value = 82.697546
$ws.Range("AW29").Value = 91.79103
$ws.Range("AW30").Value = 21.931076
$ws.Range("AW31").Value = 163.704317
$ws.Range("AW32").Value = 20.712477
$ws.Range("AW33").Value = 99.933032
$ws.Range("AW34").Value = 8.830822
$ws.Range("AW35").Value = 127.978831
$ws.Range("AW36").Value = 135.764572
$ws.Range("AK37").Value = 28.77706
$ws.Range("AW38").Value = 22.99647
$ws.Range("AK39").Value = 28.776956
$ws.Range("AW40").Value = 1.810718
$ws.Range("AN41").Value = 2.734363
$ws.Range("AN42").Value = 2.803264
$ws.Range("AW43").Value = 85.015833
$ws.Range("AW44").Value = 99.932419
$ws.Range("AW45").Value = 51.04919
$ws.Range("AW46").Value = 51.047975
$ws.Range("AK47").Value = 50.826574
$ws.Range("AW48").Value = 127.979815
$ws.Range("AW49").Value = 128.673322
$ws.Range("AW50").Value = 35.815938
$ws.Range("AW51").Value = 141.888808
$ws.Range("AW52").Value = 119.646343
$ws.Range("AW53").Value = 166.750995
$ws.Range("AW54").Value = 156.799815
$ws.Range("AQ55").Value = 0.836944
$ws.Range("AW56").Value = 128.672488
$ws.Range("AW57").Value = 23.000579
$ws.Range("AW58").Value = 2.62022
$ws.Range("AW59").Value = 64.674063
$ws.Range("AW60").Value = 131.928171
$ws.Range("AW61").Value = 92.729317
$ws.Range("AK62").Value = 23.802789
$ws.Range("AW63").Value = 12.72919
$ws.Range("AW64").Value = 168.044132
$ws.Range("AW65").Value = 170.937025
$ws.Range("AW66").Value = 128.688264
$ws.Range("AW67").Value = 85.789792
$ws.Range("AW68").Value = 21.826389
$ws.Range("AW69").Value = 28.82235
$ws.Range("AK70").Value = 28.777002
$ws.Range("AN71").Value = 2.80419
$ws.Range("AN72").Value = 2.73978
$ws.Range("AW73").Value = 142.759699
$ws.Range("AW74").Value = 79.785914
$ws.Range("AW75").Value = 14.712905
$ws.Range("AW76").Value = 83.99456
$ws.Range("AW77").Value = 1.801065
$ws.Range("AW78").Value = 156.803206
$ws.Range("AW79").Value = 51.04897
$ws.Range("AW80").Value = 35.964329
$ws.Range("AW81").Value = 30.659086
$ws.Range("AW82").Value = 20.711134
$ws.Range("AW83").Value = 135.82853
$ws.Range("AW84").Value = 166.703275
$ws.Range("AW85").Value = 114.67294
$ws.Range("AW86").Value = 83.992847
$ws.Range("AW87").Value = 64.66691
$ws.Range("AQ88").Value = 46.784792
$ws.Range("AW89").Value = 120.675451
$ws.Range("AW90").Value = 161.608171
$ws.Range("AW91").Value = 71.0486
$ws.Range("AK92").Value = 23.800347
$ws.Range("AW93").Value = 128.671713
$ws.Range("AK94").Value = 28.77691
$ws.Range("AW95").Value = 168.044525
$ws.Range("AW96").Value = 79.65919
$ws.Range("AW97").Value = 1.814294
$ws.Range("AW98").Value = 97.728368
$ws.Range("AW99").Value = 89.804225
$ws.Range("AW100").Value = 163.794132
$ws.Range("AW101").Value = 119.791123
$ws.Range("AW102").Value = 59.049873
$ws.Range("AK103").Value = 23.802813
$ws.Range("AW104").Value = 0.82919
$ws.Range("AW105").Value = 118.227998
$ws.Range("AW106").Value = 126.871516
$ws.Range("AW107").Value = 79.715139
$ws.Range("AK108").Value = 28.77713
$ws.Range("AW109").Value = 28.803831
$ws.Range("AK110").Value = 28.776944
$ws.Range("AW111").Value = 1.799745
$ws.Range("AW112").Value = 170.939988
$ws.Range("AW113").Value = 156.750868
$ws.Range("AW114").Value = 166.701412
$ws.Range("AW115").Value = 59.048275
$ws.Range("AW116").Value = 30.649722
$ws.Range("AW117").Value = 12.828553
$ws.Range("AW118").Value = 145.670081
$ws.Range("AW119").Value = 119.811134
$ws.Range("AW120").Value = 97.681574
$ws.Range("AW121").Value = 83.990694
$ws.Range("AK122").Value = 28.777072
$ws.Range("AK123").Value = 51.626921
$ws.Range("AW124").Value = 161.608785
$ws.Range("AW125").Value = 1.801169
$ws.Range("AN126").Value = 2.738356
$ws.Range("AW127").Value = 156.70338
$ws.Range("AQ128").Value = 0.8361
$ws.Range("AK129").Value = 28.776956
$ws.Range("AW130").Value = 70.910949
$ws.Range("AW131").Value = 65.062361
$ws.Range("AW132").Value = 174.753403
$ws.Range("AW133").Value = 100.593762
$ws.Range("AW134").Value = 59.047731
$ws.Range("AQ135").Value = 18.716019
$ws.Range("AK136").Value = 28.777049
$ws.Range("AW137").Value = 156.754074
$ws.Range("AW138").Value = 161.61044
$ws.Range("AW139").Value = 97.791065
$ws.Range("AK140").Value = 28.776991
$ws.Range("AW141").Value = 131.927743
$ws.Range("AW142").Value = 128.672211
$ws.Range("AW143").Value = 91.79206
$ws.Range("AK144").Value = 28.777234
$ws.Range("AW145").Value = 63.924236
$ws.Range("AW146").Value = 63.724294
$ws.Range("AQ147").Value = 46.784606
$ws.Range("AW148").Value = 85.622431
$ws.Range("AW149").Value = 75.687662
$ws.Range("AW150").Value = 105.757836
$ws.Range("AK151").Value = 28.777095
$ws.Range("AW152").Value = 59.047037
$ws.Range("AW153").Value = 2.62088
$ws.Range("AW154").Value = 127.979618
$ws.Range("AW155").Value = 83.996308
$ws.Range("AK156").Value = 23.800347
$ws.Range("AW157").Value = 7.643831
$ws.Range("AW158").Value = 118.002431
$ws.Range("AW159").Value = 33.748877
$ws.Range("AW160").Value = 1.79809
$ws.Range("AW161").Value = 166.704213
$ws.Range("AW162").Value = 143.054514
$ws.Range("AW163").Value = 82.697338
$ws.Range("AK164").Value = 28.776921
$ws.Range("AW165").Value = 118.007373
$ws.Range("AW166").Value = 170.940382
$ws.Range("AW167").Value = 92.708449
$ws.Range("AW168").Value = 163.703762
$ws.Range("AW169").Value = 154.810602
$ws.Range("AW170").Value = 92.729201
$ws.Range("AW171").Value = 21.826227
$ws.Range("AK172").Value = 50.826563
$ws.Range("AW173").Value = 163.703403
$ws.Range("AW174").Value = 119.806944
$ws.Range("AW175").Value = 65.061644
$ws.Range("AW176").Value = 97.621354
$ws.Range("AW177").Value = 8.774699
$ws.Range("AW178").Value = 170.947951
$ws.Range("AW179").Value = 105.770023
$ws.Range("AW180").Value = 83.993819
$ws.Range("AW181").Value = 65.065185
$ws.Range("AW182").Value = 119.810764
$ws.Range("AK183").Value = 28.777072
$ws.Range("AQ184").Value = 0.836343
$ws.Range("AW185").Value = 163.79441
$ws.Range("AW186").Value = 178.015243
$ws.Range("AW187").Value = 83.999664
$ws.Range("AW188").Value = 57.78706
$ws.Range("AW189").Value = 151.030949
$ws.Range("AW190").Value = 161.609352
$ws.Range("AW191").Value = 89.712905
$ws.Range("AW192").Value = 145.666921
$ws.Range("AW193").Value = 93.636343
$ws.Range("AK194").Value = 23.802801
$ws.Range("AW195").Value = 166.677546
$ws.Range("AW196").Value = 163.793866
$ws.Range("AW197").Value = 162.716088
$ws.Range("AW198").Value = 174.698669
$ws.Range("AW199").Value = 176.839063
$ws.Range("AW200").Value = 166.752245
$ws.Range("AW201").Value = 128.673935
$ws.Range("AW202").Value = 34.016944
$ws.Range("AK203").Value = 50.826574
$ws.Range("AW204").Value = 128.69191
$ws.Range("AW205").Value = 175.794387
$ws.Range("AW206").Value = 166.69331
$ws.Range("AW207").Value = 175.790255
$ws.Range("AW208").Value = 76.74456
$ws.Range("AK209").Value = 28.777025
$ws.Range("AK210").Value = 28.777014
$ws.Range("AW211").Value = 8.798854
$ws.Range("AW212").Value = 1.817002
$ws.Range("AW213").Value = 83.999016
$ws.Range("AW214").Value = 135.764491
$ws.Range("AW215").Value = 153.005799
$ws.Range("AW216").Value = 92.775405
$ws.Range("AW217").Value = 85.015602
$ws.Range("AW218").Value = 65.063495
$ws.Range("AW219").Value = 15.992361
$ws.Range("AW220").Value = 127.98015
$ws.Range("AW221").Value = 83.858611
$ws.Range("AQ222").Value = 0.835822
$ws.Range("AW223").Value = 2.621377
$ws.Range("AW224").Value = 118.006146
$ws.Range("AW225").Value = 163.678646
$ws.Range("AW226").Value = 128.67272
$ws.Range("AW227").Value = 92.706794
$ws.Range("AW228").Value = 35.961377
$ws.Range("AW229").Value = 64.766262
$ws.Range("AW230").Value = 59.046157
$ws.Range("AK231").Value = 28.776968
$ws.Range("AK232").Value = 51.626933
$ws.Range("AW233").Value = 128.673125
$ws.Range("AK234").Value = 28.777049
$ws.Range("AW235").Value = 99.933657
$ws.Range("AW236").Value = 148.826319
$ws.Range("AK237").Value = 28.770046
$ws.Range("AK238").Value = 28.776921
$ws.Range("AW239").Value = 127.940729
$ws.Range("AW240").Value = 59.055903
$ws.Range("AW241").Value = 59.05375
$ws.Range("AW242").Value = 33.638299
$ws.Range("AK243").Value = 28.777002
$ws.Range("AW244").Value = 125.884988
$ws.Range("AW245").Value = 50.044468
$ws.Range("AW246").Value = 35.81669
$ws.Range("AW247").Value = 156.803426
$ws.Range("AW248").Value = 166.751505
$ws.Range("AW249").Value = 145.667269
$ws.Range("AW250").Value = 21.826123
$ws.Range("AW251").Value = 35.8175
$ws.Range("AW252").Value = 1.800891
$ws.Range("AW253").Value = 145.667813
$ws.Range("AW254").Value = 65.057894
$ws.Range("AW255").Value = 16.804734
$ws.Range("AW256").Value = 155.69662
$ws.Range("AW257").Value = 70.912141
$ws.Range("AW258").Value = 22.99875
$ws.Range("AW259").Value = 141.888611
$ws.Range("AW260").Value = 21.931806
$ws.Range("AW261").Value = 12.727813
$ws.Range("AW262").Value = 1.817662
$ws.Range("AW263").Value = 51.047894
$ws.Range("AQ264").Value = 11.627789
$ws.Range("AW265").Value = 143.054884
$ws.Range("AW266").Value = 154.811551
$ws.Range("AW267").Value = 118.007211
$ws.Range("AW268").Value = 97.735278
$ws.Range("AW269").Value = 39.062697
$ws.Range("AW270").Value = 16.798391
$ws.Range("AW271").Value = 142.01728
$ws.Range("AK272").Value = 50.826586
$ws.Range("AW273").Value = 156.016481
$ws.Range("AW274").Value = 100.624375
$ws.Range("AW275").Value = 161.607384
$ws.Range("AW276").Value = 84.000231
$ws.Range("AW277").Value = 128.675058
$ws.Range("AW278").Value = 1.66912
$ws.Range("AW279").Value = 153.006262
$ws.Range("AW280").Value = 83.993102
$ws.Range("AW281").Value = 58.804572
$ws.Range("AK282").Value = 28.777176
$ws.Range("AW283").Value = 93.766042
$ws.Range("AW284").Value = 168.045984
$ws.Range("AW285").Value = 57.718194
$ws.Range("AK286").Value = 28.776933
$ws.Range("AW287").Value = 12.828403
$ws.Range("AQ288").Value = 0.835336
$ws.Range("AW289").Value = 57.716759
$ws.Range("AW290").Value = 57.71853
$ws.Range("AK291").Value = 51.626933
$ws.Range("AW292").Value = 161.607766
$ws.Range("AW293").Value = 114.669456
$ws.Range("AW294").Value = 33.749826
$ws.Range("AW295").Value = 114.672245
$ws.Range("AW296").Value = 127.977488
$ws.Range("AK297").Value = 23.800336
$ws.Range("AW298").Value = 170.940833
$ws.Range("AW299").Value = 65.057037
$ws.Range("AW300").Value = 28.767118
$ws.Range("AW301").Value = 16.802373
$ws.Range("AW302").Value = 149.714884
$ws.Range("AW303").Value = 104.7739
$ws.Range("AN304").Value = 2.803275
$ws.Range("AN305").Value = 2.803322
$ws.Range("AW306").Value = 128.674387
$ws.Range("AW307").Value = 97.681273
$ws.Range("AW308").Value = 43.691991
$ws.Range("AW309").Value = 22.994965
$ws.Range("AW310").Value = 12.728241
$ws.Range("AW311").Value = 107.698866
$ws.Range("AK312").Value = 28.777176
$ws.Range("AW313").Value = 99.930301
$ws.Range("AW314").Value = 33.74934
$ws.Range("AW315").Value = 175.799213
$ws.Range("AW316").Value = 85.998183
$ws.Range("AW317").Value = 65.064398
$ws.Range("AW318").Value = 14.712419
$ws.Range("AW319").Value = 133.037558
$ws.Range("AK320").Value = 28.777106
$ws.Range("AK321").Value = 28.777153
$ws.Range("AQ322").Value = 11.633796
$ws.Range("AW323").Value = 89.75794
$ws.Range("AW324").Value = 68.696458
$ws.Range("AW325").Value = 33.750451
$ws.Range("AW326").Value = 14.713646
$ws.Range("AW327").Value = 142.0175
$ws.Range("AW328").Value = 70.91169
$ws.Range("AW329").Value = 59.045394
$ws.Range("AW330").Value = 65.117905
$ws.Range("AK331").Value = 28.777222
$ws.Range("AW332").Value = 0.83265
$ws.Range("AW333").Value = 166.702118
$ws.Range("AW334").Value = 168.044294
$ws.Range("AK335").Value = 28.777118
$ws.Range("AN336").Value = 2.803252
$ws.Range("AW337").Value = 174.692963
$ws.Range("AW338").Value = 65.065613
$ws.Range("AW339").Value = 145.668333
$ws.Range("AW340").Value = 35.962882
$ws.Range("AW341").Value = 33.997917
$ws.Range("AW342").Value = 161.607118
$ws.Range("AW343").Value = 133.034873
$ws.Range("AW344").Value = 51.635069
$ws.Range("AW345").Value = 135.827431
$ws.Range("AW346").Value = 83.994456
$ws.Range("AW347").Value = 170.83287
$ws.Range("AW348").Value = 166.738727
$ws.Range("AW349").Value = 14.719039
$ws.Range("AW350").Value = 168.045093
$ws.Range("AW351").Value = 145.668854
$ws.Range("AW352").Value = 175.790613
$ws.Range("AW353").Value = 23.726991
$ws.Range("AW354").Value = 161.605648
$ws.Range("AW355").Value = 128.671528
$ws.Range("AW356").Value = 119.804155
$ws.Range("AW357").Value = 35.962199
$ws.Range("AN358").Value = 2.803299
$ws.Range("AN359").Value = 2.804167
$ws.Range("AW360").Value = 82.697685
$ws.Range("AW361").Value = 82.697222
$ws.Range("AQ362").Value = 0.835556
$ws.Range("AW363").Value = 14.712627
$ws.Range("AW364").Value = 1.815266
